$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in newly-collected experiment results (rows 44-46, cols G-J) ---
$ws.Range("I44").Value = 3.8729328155517502
$ws.Range("I44").NumberFormat = "#,##0"
$ws.Range("J44").Value = 1.05591082663456

$ws.Range("G45").Value = 0.654283475875854
$ws.Range("H45").Value = 0.060748244857041098
$ws.Range("I45").Value = 2.9005457878112701
$ws.Range("J45").Value = 0.16064807726299399

$ws.Range("G46").Value = 912.27280468940705

# --- Drop the stray leftover value beneath the table ---
$ws.Range("H48").ClearContents()

# --- Add a new underline-styled spacer cell (matching F49 / G52) ---
$ws.Range("F49").Copy()
$ws.Range("F51").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Move the saved selection / scroll position back to the top ---
$ws.Range("C5").Select() | Out-Null

# --- Set page layout (paper size + orientation) used for printing ---
$ws.PageSetup.PaperSize = 9        # xlPaperA4
$ws.PageSetup.Orientation = 1      # xlPortrait
